$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.771609
$ws.Range("H2").Value = 26.314827
$ws.Range("I2").Value = 0.2200338127677125
$ws.Range("J2").Value = 0.2200338127677125
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 90.76522116705632
$ws.Range("R2").Value = 816.886990503507
$ws.Range("S2").Value = 0.05060675328694966
$ws.Range("T2").Value = 0.05060675328694967

# Row 3
$ws.Range("G3").Value = 8.771609
$ws.Range("H3").Value = 26.314827
$ws.Range("I3").Value = 0.2200338127677125
$ws.Range("J3").Value = 0.2200338127677125
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 268.1468795718276
$ws.Range("R3").Value = 2413.321916146449
$ws.Range("S3").Value = 0.1495070777625362
$ws.Range("T3").Value = 0.1495070777625362

# Row 4
$ws.Range("G4").Value = 8.771609
$ws.Range("H4").Value = 26.314827
$ws.Range("I4").Value = 0.2200338127677125
$ws.Range("J4").Value = 0.2200338127677125
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 35.72727805806133
$ws.Range("R4").Value = 321.545502522552
$ws.Range("S4").Value = 0.01991998171822664
$ws.Range("T4").Value = 0.01991998171822664

# Row 5
$ws.Range("I5").Value = 0.583164828467109
$ws.Range("J5").Value = 0.583164828467109
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 240.5588666890231
$ws.Range("R5").Value = 2165.029800201208
$ws.Range("S5").Value = 0.134125197525968
$ws.Range("T5").Value = 0.134125197525968

# Row 6
$ws.Range("I6").Value = 0.583164828467109
$ws.Range("J6").Value = 0.583164828467109
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("S6").Value = 0.3962448691922222
$ws.Range("T6").Value = 0.3962448691922222

# Row 7
$ws.Range("I7").Value = 0.583164828467109
$ws.Range("J7").Value = 0.583164828467109
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 94.68950121007646
$ws.Range("R7").Value = 852.2055108906881
$ws.Range("S7").Value = 0.05279476174891879
$ws.Range("T7").Value = 0.0527947617489188

# Row 8
$ws.Range("G8").Value = 7.845451333333334
$ws.Range("H8").Value = 23.536354
$ws.Range("I8").Value = 0.1968013587651783
$ws.Range("J8").Value = 0.1968013587651783
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 81.18169943796822
$ws.Range("R8").Value = 730.635294941714
$ws.Range("S8").Value = 0.04526339694926783
$ws.Range("T8").Value = 0.04526339694926785

# Row 9
$ws.Range("G9").Value = 7.845451333333334
$ws.Range("H9").Value = 23.536354
$ws.Range("I9").Value = 0.1968013587651783
$ws.Range("J9").Value = 0.1968013587651783
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 239.8343672028665
$ws.Range("R9").Value = 2158.509304825798
$ws.Range("S9").Value = 0.1337212480144589
$ws.Range("T9").Value = 0.1337212480144589

# Row 10
$ws.Range("G10").Value = 7.845451333333334
$ws.Range("H10").Value = 23.536354
$ws.Range("I10").Value = 0.1968013587651783
$ws.Range("J10").Value = 0.1968013587651783
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 31.95498354714489
$ws.Range("R10").Value = 287.594851924304
$ws.Range("S10").Value = 0.01781671380145157
$ws.Range("T10").Value = 0.01781671380145158
